$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new headers I1 and J1, copying the style used by the existing header row (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF) for rows 2..33
$data = @{
    2  = @(7, 7)
    3  = @(6, 6)
    4  = @(6, 6)
    5  = @(7, 7)
    6  = @(10, 10)
    7  = @(8, 8)
    8  = @(7, 7)
    9  = @(9, 9)
    10 = @(7, 7)
    11 = @(7, 7)
    12 = @(6, 7)
    13 = @(7, 7)
    14 = @(7, 7)
    15 = @(7, 7)
    16 = @(1, 6)
    17 = @(1, 6)
    18 = @(1, 5)
    19 = @(1, 6)
    20 = @(1, 5)
    21 = @(1, 6)
    22 = @(1, 6)
    23 = @(1, 4)
    24 = @(1, 6)
    25 = @(1, 4)
    26 = @(1, 5)
    27 = @(1, 4)
    28 = @(1, 4)
    29 = @(1, 3)
    30 = @(1, 2)
    31 = @(1, 3)
    32 = @(1, 2)
    33 = @(1, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
